$wb = $excel.ActiveWorkbook

# zh-cn sheet: rows 2 and 4 share the same "Correspond Handoff/Handback
# Datetime" values (15d158d6... and 63aea1f3... rows), so both must be
# updated together.
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-18 03:47:54"
$wsZhCn.Range("H2").Value = "2016-03-18 03:48:13"
$wsZhCn.Range("E4").Value = "2016-03-18 03:47:54"
$wsZhCn.Range("H4").Value = "2016-03-18 03:48:13"

# de-de sheet: same pairing of rows 2 and 4.
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-18 03:47:58"
$wsDeDe.Range("H2").Value = "2016-03-18 03:48:19"
$wsDeDe.Range("E4").Value = "2016-03-18 03:47:58"
$wsDeDe.Range("H4").Value = "2016-03-18 03:48:19"
